$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D1").Value = "EmpID"
$headerRange = $ws.Range("A1:D1")
$font = $headerRange.Font
$font.Size = 11
$font.Bold = $true
$ws.Rows(1).RowHeight = 15
$headerRange.Select() | Out-Null
